# CWND -> CongestionWindow for CUBIC
#
# Slide 9, shape "TextBox 11" (id 12) holds the axis label text "CWND".
# Rename it to "CongestionWindow" and rotate/resize/reposition the text
# box so the (now much longer) label still reads bottom-to-top along the
# vertical axis of the plot.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Locate the textbox by name rather than a hard-coded index, in case
# shape ordering ever shifts.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 11") {
        $shp = $candidate
        break
    }
}

$shp.TextFrame.TextRange.Text = "CongestionWindow"

# Reposition/resize/rotate the box (EMU values below, expressed in points
# as the object model expects, chosen so the float32 Shape.Left/Top/
# Width/Height round-trip to the exact target EMU offsets):
#   off  x=603888  y=6473916
#   ext  cx=4275529 cy=707886
#   rot  16200000 (270 degrees)
$shp.Left = 47.55023956298828
$shp.Top = 509.7571716308594
$shp.Width = 336.6558532714844
$shp.Height = 55.73905563354492
$shp.Rotation = 270
